$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1146.3636
$ws.Range("J28").Value = 1342.2858
$ws.Range("L28").Value = 1342.2858
$ws.Range("N28").Value = -2312.2858

$ws.Range("H40").Value = 35837.168
$ws.Range("J40").Value = 36666
$ws.Range("L40").Value = 36666
$ws.Range("N40").Value = -37016

$ws.Range("H86").Value = 111114910
$ws.Range("J86").Value = 5499.5
$ws.Range("L86").Value = 5499.5
$ws.Range("N86").Value = -7745.5

$ws.Range("H89").Value = 111114910
$ws.Range("J89").Value = 5499.5
$ws.Range("L89").Value = 27497.5
$ws.Range("N89").Value = -38729.5

$ws.Range("H106").Value = 3450
$ws.Range("I106").Value = 3322.9092
$ws.Range("K106").Value = 3322.9092
$ws.Range("M106").Value = -2691.9092

$ws.Range("H125").Value = 2958.9375
$ws.Range("I125").Value = 2330.111
$ws.Range("J125").Value = 3767.4285
$ws.Range("K125").Value = 20970.999
$ws.Range("L125").Value = 33906.8565
$ws.Range("M125").Value = -18510.999
$ws.Range("N125").Value = -38826.8565

$ws.Range("H137").Value = 10099789
$ws.Range("I137").Value = 589502.7
$ws.Range("J137").Value = 19610076
$ws.Range("K137").Value = 1768508.1
$ws.Range("L137").Value = 58830228
$ws.Range("M137").Value = -1765958.1
$ws.Range("N137").Value = -58835328

$ws.Range("H138").Value = 1868.63
$ws.Range("I138").Value = 1105.4166
$ws.Range("J138").Value = 1972.7046
$ws.Range("K138").Value = 3316.2498
$ws.Range("L138").Value = 5918.1138
$ws.Range("M138").Value = 1823.7502
$ws.Range("N138").Value = -16198.1138

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 10000
$ws.Range("J43").Value = 10000
$ws.Range("L43").Value = 10000
$ws.Range("N43").Value = -10626

$ws.Range("H74").Value = 17858930
$ws.Range("I74").Value = 25001332
$ws.Range("K74").Value = 25001332
$ws.Range("M74").Value = -25000458

$ws.Range("H77").Value = 17858930
$ws.Range("I77").Value = 25001332
$ws.Range("K77").Value = 125006660
$ws.Range("M77").Value = -125002292

$ws.Range("H80").Value = 67994.5
$ws.Range("J80").Value = 67994.5
$ws.Range("L80").Value = 67994.5
$ws.Range("N80").Value = -69990.5

$ws.Range("H83").Value = 67994.5
$ws.Range("J83").Value = 67994.5
$ws.Range("L83").Value = 203983.5
$ws.Range("N83").Value = -213967.5

$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

$ws.Range("H97").Value = 565.3684
$ws.Range("I97").Value = 594.4667
$ws.Range("K97").Value = 594.4667
$ws.Range("M97").Value = -98.46669999999995

$ws.Range("H137").Value = 84999.664
$ws.Range("J137").Value = 72500
$ws.Range("L137").Value = 72500
$ws.Range("N137").Value = -82700

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 915.6667
$ws.Range("I16").Value = 915.6667
$ws.Range("K16").Value = 915.6667
$ws.Range("M16").Value = -628.6667

$ws.Range("H51").Value = 45000
$ws.Range("J51").Value = 45000
$ws.Range("L51").Value = 45000
$ws.Range("N51").Value = -46472

$ws.Range("H58").Value = 2108.743
$ws.Range("I58").Value = 1974.8518
$ws.Range("K58").Value = 1974.8518
$ws.Range("M58").Value = -1771.8518

$ws.Range("H61").Value = 45000
$ws.Range("J61").Value = 45000
$ws.Range("L61").Value = 45000
$ws.Range("N61").Value = -45696

$ws.Range("H107").Value = 810.93335
$ws.Range("J107").Value = 1048.625
$ws.Range("L107").Value = 1048.625
$ws.Range("N107").Value = -4888.625

$ws.Range("H113").Value = 915.6667
$ws.Range("I113").Value = 915.6667
$ws.Range("K113").Value = 915.6667
$ws.Range("M113").Value = 1254.3333

$ws.Range("H132").Value = 35090200
$ws.Range("I132").Value = 39217804
$ws.Range("K132").Value = 117653412
$ws.Range("M132").Value = -117650882

$ws.Range("H134").Value = 2193.3125
$ws.Range("I134").Value = 2167.2903
$ws.Range("K134").Value = 6501.8709
$ws.Range("M134").Value = -3966.8709

$ws.Range("H136").Value = 2108.743
$ws.Range("I136").Value = 1974.8518
$ws.Range("K136").Value = 5924.555399999999
$ws.Range("M136").Value = -3374.555399999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 3810.4443
$ws.Range("J39").Value = 3978.4285
$ws.Range("L39").Value = 11935.2855
$ws.Range("N39").Value = -12523.2855

$ws.Range("H75").Value = 3242.5
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 3242.5
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 9727.5
$ws.Range("N75").Value = -11723.5
$ws.Range("M75").ClearContents()

$ws.Range("H78").Value = 3242.5
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 3242.5
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 29182.5
$ws.Range("N78").Value = -39166.5
$ws.Range("M78").ClearContents()

$ws.Range("H128").Value = 154499.5
$ws.Range("I128").Value = 154499.5
$ws.Range("K128").Value = 463498.5
$ws.Range("M128").Value = -458518.5

$ws.Range("H134").Value = 4404.8125
$ws.Range("I134").Value = 4404.8125
$ws.Range("K134").Value = 13214.4375
$ws.Range("M134").Value = -8144.4375

$ws.Range("H139").Value = 2147.5757
$ws.Range("I139").Value = 1360.6154
$ws.Range("J139").Value = 5070.5713
$ws.Range("K139").Value = 4081.8462
$ws.Range("L139").Value = 15211.7139
$ws.Range("M139").Value = 1058.1538
$ws.Range("N139").Value = -25491.7139

$ws.Range("H140").Value = 3984
$ws.Range("I140").Value = 3057.3076
$ws.Range("K140").Value = 9171.9228
$ws.Range("M140").Value = -3991.9228

$ws.Range("H141").Value = 4047.158
$ws.Range("I141").Value = 3118.5625
$ws.Range("J141").Value = 8999.667
$ws.Range("K141").Value = 9355.6875
$ws.Range("L141").Value = 26999.001
$ws.Range("M141").Value = -4175.6875
$ws.Range("N141").Value = -37359.001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 307466.6
$ws.Range("I122").Value = 478542.84
$ws.Range("J122").Value = 8083.1665
$ws.Range("K122").Value = 1435628.52
$ws.Range("L122").Value = 24249.4995
$ws.Range("M122").Value = -1433178.52
$ws.Range("N122").Value = -29149.4995

$ws.Range("H135").Value = 124999.5
$ws.Range("J135").Value = 124999.5
$ws.Range("L135").Value = 124999.5
$ws.Range("N135").Value = -135139.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3043.84
$ws.Range("I61").Value = 3244.6086
$ws.Range("K61").Value = 3244.6086
$ws.Range("M61").Value = -3042.6086

$ws.Range("H113").Value = 3043.84
$ws.Range("I113").Value = 3244.6086
$ws.Range("K113").Value = 3244.6086
$ws.Range("M113").Value = -1074.6086

$ws.Range("H139").Value = 140000
$ws.Range("J139").Value = 140000
$ws.Range("L139").Value = 140000
$ws.Range("N139").Value = -150280

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 10103937
$ws.Range("I132").Value = 16668697
$ws.Range("K132").Value = 50006091
$ws.Range("M132").Value = -50003561

$ws.Range("H136").Value = 3703.5625
$ws.Range("I136").Value = 2563.0212
$ws.Range("K136").Value = 7689.0636
$ws.Range("M136").Value = -5139.0636
